$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$st = $ws.Range("D2").Style
$ws.Range("D2").Value = '''61.152.92'
$ws.Range("D2").Style = $st
$ws.Range("E2").Value = '  +6.96%  '

$st = $ws.Range("D3").Style
$ws.Range("D3").Value = '''3.313.26'
$ws.Range("D3").Style = $st
$ws.Range("E3").Value = '  +2.02%  '

$st = $ws.Range("D4").Style
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = $st
$ws.Range("E4").Value = '  +0.23%  '

$st = $ws.Range("D5").Style
$ws.Range("D5").Value = '''410.38'
$ws.Range("D5").Style = $st
$ws.Range("E5").Value = '  +3.75%  '

$st = $ws.Range("D6").Style
$ws.Range("D6").Value = '''112.28'
$ws.Range("D6").Style = $st
$ws.Range("E6").Value = '  +4.01%  '

$st = $ws.Range("D7").Style
$ws.Range("D7").Value = '''3.308.44'
$ws.Range("D7").Style = $st
$ws.Range("E7").Value = '  +2.00%  '

$st = $ws.Range("D8").Style
$ws.Range("D8").Value = '''0.566'
$ws.Range("D8").Style = $st
$ws.Range("E8").Value = '  -3.92%  '

$ws.Range("E9").Value = '  +0.03%  '

$st = $ws.Range("D10").Style
$ws.Range("D10").Value = '''0.619'
$ws.Range("D10").Style = $st
$ws.Range("E10").Value = '  -1.10%  '

$st = $ws.Range("D11").Style
$ws.Range("D11").Value = '''0.115'
$ws.Range("D11").Style = $st
$ws.Range("E11").Value = '  +16.91%  '

$st = $ws.Range("D12").Style
$ws.Range("D12").Value = '''38.73'
$ws.Range("D12").Style = $st
$ws.Range("E12").Value = '  -1.09%  '

$ws.Range("E13").Value = '  -0.17%  '

$st = $ws.Range("D14").Style
$ws.Range("D14").Value = '''3.847.95'
$ws.Range("D14").Style = $st
$ws.Range("E14").Value = '  +2.36%  '

$st = $ws.Range("D15").Style
$ws.Range("D15").Value = '''8.14'
$ws.Range("D15").Style = $st
$ws.Range("E15").Value = '  -0.88%  '

$st = $ws.Range("D16").Style
$ws.Range("D16").Value = '''18.95'
$ws.Range("D16").Style = $st
$ws.Range("E16").Value = '  -0.79%  '

$st = $ws.Range("D17").Style
$ws.Range("D17").Value = '''3.338.04'
$ws.Range("D17").Style = $st
$ws.Range("E17").Value = '  +3.56%  '

$st = $ws.Range("D18").Style
$ws.Range("D18").Value = '''61.119.46'
$ws.Range("D18").Style = $st
$ws.Range("E18").Value = '  +7.43%  '

$st = $ws.Range("D19").Style
$ws.Range("D19").Value = '''0.983'
$ws.Range("D19").Style = $st
$ws.Range("E19").Value = '  -4.30%  '

$st = $ws.Range("D20").Style
$ws.Range("D20").Value = '''10.53'
$ws.Range("D20").Style = $st
$ws.Range("E20").Value = '  -2.50%  '

$st = $ws.Range("D21").Style
$ws.Range("D21").Value = '''0.0000114'
$ws.Range("D21").Style = $st
$ws.Range("E21").Value = '  +1.22%  '

$st = $ws.Range("D22").Style
$ws.Range("D22").Value = '''3.22'
$ws.Range("D22").Style = $st
$ws.Range("E22").Value = '  -3.98%  '

$ws.Range("E23").Value = '  -4.53%  '

$st = $ws.Range("D24").Style
$ws.Range("D24").Value = '''294.64'
$ws.Range("D24").Style = $st
$ws.Range("E24").Value = '  -0.44%  '

$st = $ws.Range("D25").Style
$ws.Range("D25").Value = '''73.09'
$ws.Range("D25").Style = $st
$ws.Range("E25").Value = '  -1.61%  '

$st = $ws.Range("D26").Style
$ws.Range("D26").Value = '''3.07'
$ws.Range("D26").Style = $st
$ws.Range("E26").Value = '  -2.62%  '

$st = $ws.Range("D27").Style
$ws.Range("D27").Value = '''29.10'
$ws.Range("D27").Style = $st
$ws.Range("E27").Value = '  +4.30%  '

$st = $ws.Range("D28").Style
$ws.Range("D28").Value = '''4.56'
$ws.Range("D28").Style = $st
$ws.Range("E28").Value = '  +4.48%  '

$ws.Range("E29").Value = '  +2.71%  '

$st = $ws.Range("D30").Style
$ws.Range("D30").Value = '''7.33'
$ws.Range("D30").Style = $st
$ws.Range("E30").Value = '  +1.10%  '

$st = $ws.Range("D31").Style
$ws.Range("D31").Value = '''7.40'
$ws.Range("D31").Style = $st
$ws.Range("E31").Value = '  -2.84%  '

$st = $ws.Range("D32").Style
$ws.Range("D32").Value = '''1.00'
$ws.Range("D32").Style = $st
$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("E33").Value = '  +0.61%  '

$st = $ws.Range("D34").Style
$ws.Range("D34").Value = '''11.10'
$ws.Range("D34").Style = $st
$ws.Range("E34").Value = '  -2.38%  '

$st = $ws.Range("D35").Style
$ws.Range("D35").Value = '''2.44'
$ws.Range("D35").Style = $st
$ws.Range("E35").Value = '  +14.77%  '

$st = $ws.Range("D36").Style
$ws.Range("D36").Value = '''39.38'
$ws.Range("D36").Style = $st
$ws.Range("E36").Value = '  +0.74%  '

$st = $ws.Range("D37").Style
$ws.Range("D37").Value = '''0.0478'
$ws.Range("D37").Style = $st
$ws.Range("E37").Value = '  -0.59%  '

$st = $ws.Range("D38").Style
$ws.Range("D38").Value = '''52.50'
$ws.Range("D38").Style = $st
$ws.Range("E38").Value = '  +1.78%  '

$st = $ws.Range("D39").Style
$ws.Range("D39").Value = '''0.996'
$ws.Range("D39").Style = $st
$ws.Range("E39").Value = '  -0.15%  '

$st = $ws.Range("D40").Style
$ws.Range("D40").Value = '''3.02'
$ws.Range("D40").Style = $st
$ws.Range("E40").Value = '  +3.31%  '

$st = $ws.Range("D41").Style
$ws.Range("D41").Value = '''3.27'
$ws.Range("D41").Style = $st
$ws.Range("E41").Value = '  -6.37%  '

$st = $ws.Range("D42").Style
$ws.Range("D42").Value = '''134.80'
$ws.Range("D42").Style = $st
$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("E43").Value = '  -1.72%  '

$st = $ws.Range("D44").Style
$ws.Range("D44").Value = '''1.87'
$ws.Range("D44").Style = $st
$ws.Range("E44").Value = '  -0.56%  '

$st = $ws.Range("D45").Style
$ws.Range("D45").Value = '''0.282'
$ws.Range("D45").Style = $st
$ws.Range("E45").Value = '  +0.35%  '

$st = $ws.Range("D46").Style
$ws.Range("D46").Value = '''16.19'
$ws.Range("D46").Style = $st
$ws.Range("E46").Value = '  -5.08%  '

$st = $ws.Range("D47").Style
$ws.Range("D47").Value = '''3.74'
$ws.Range("D47").Style = $st
$ws.Range("E47").Value = '  -5.06%  '

$ws.Range("E48").Value = '  +2.54%  '

$st = $ws.Range("D49").Style
$ws.Range("D49").Value = '''20.83'
$ws.Range("D49").Style = $st
$ws.Range("E49").Value = '  -6.09%  '

$st = $ws.Range("D50").Style
$ws.Range("D50").Value = '''2.108.02'
$ws.Range("D50").Style = $st
$ws.Range("E50").Value = '  -2.44%  '

$st = $ws.Range("D51").Style
$ws.Range("D51").Value = '''3.648.67'
$ws.Range("D51").Style = $st
$ws.Range("E51").Value = '  +2.25%  '

